$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 17 (this pushes the totals row 17->18 and footer row 18->19 down)
$ws.Rows("17:17").Insert()

# Copy formatting from the row above (row 16, the last data row) into the new row 17
$ws.Range("A16:Q16").Copy()
$ws.Range("A17:Q17").PasteSpecial(-4122)

# Fix up the Q17 style (PasteSpecial mis-mapped the last cell in the pasted range)
$ws.Range("Q16").Copy()
$ws.Range("Q17").PasteSpecial(-4122)

# Re-create merges for the new row 17 (PasteSpecial formats does not copy merge state)
$ws.Range("A17:B17").Merge()
$ws.Range("C17:G17").Merge()
$ws.Range("H17:K17").Merge()
$ws.Range("L17:M17").Merge()
$ws.Range("N17:O17").Merge()

# Populate the new item row (#11)
$ws.Range("A17").Value = 11
$ws.Range("C17").Value = "صوفي طويل جدا"
$ws.Range("H17").Value = "1:0"
$ws.Range("L17").Value = "0"
$ws.Range("N17").Value = "50.00"
$ws.Range("P17").Value = "50.0000"
$ws.Range("Q17").Value = "1:0"

# Row heights: new data row matches the other data rows; totals row below now uses 24.75
$ws.Rows("17:17").RowHeight = 25.5
$ws.Rows("18:18").RowHeight = 24.75

# Update the running total (now on row 18)
$ws.Range("P18").Value = 473.25

# Update the footer timestamp (now on row 19)
$ws.Range("A19").Value = "Monday, 28 July, 2025 11:03 AM"
